# Generate Report for Handoff
#
# The 9454ebbf-8be6-487a-86c3-a2e4b2f7c07d.md file has finished its
# handoff-generation pass: its translation status flips from
# "In Translation" to "Ready for handoff", its priority moves from
# human translation ("ht") to machine translation ("mt"), and the
# handoff timestamps are refreshed - on the Overview roll-up sheet and
# on each per-locale sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the 9454ebbf... file -------------------
# Columns: A=File Name, B=Path And Name, C=Extension, D=Publish URL,
#          E=zh-cn, F=de-de, G=Latest HO Xliff Generate Date
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-06 20:17:51"

# --- zh-cn sheet: row 3 is the 9454ebbf... file -----------------------
# Columns: A=Source File Name, B=File Extension, C=Status, D=Source Path,
#          E=Priority, ... H=Latest Handoff Datetime
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-09-06 20:17:45"

# --- de-de sheet: row 3 is the 9454ebbf... file -----------------------
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-09-06 20:17:51"

# --- Cosmetic: the longer "Ready for handoff" status text widens the
# columns that now display it (autofit side effect in the source
# workbook). Approximate the resulting column widths as closely as the
# host's width quantization allows.
$overview.Range("E1").ColumnWidth = 16.25
$overview.Range("F1").ColumnWidth = 16.25
$zhcn.Range("C1").ColumnWidth = 16.25
$dede.Range("C1").ColumnWidth = 16.25
